# Update "想去人数" (number of people interested) values on two sheets:
# "展览" (sheet1) and "全部类型" (sheet4), reflecting the refreshed
# generated output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 13
$wsExhibit.Range("F8").Value = 7791
$wsExhibit.Range("F11").Value = 6772
$wsExhibit.Range("F14").Value = 4808
$wsExhibit.Range("F15").Value = 5172
$wsExhibit.Range("F23").Value = 184
$wsExhibit.Range("F25").Value = 90
$wsExhibit.Range("F26").Value = 8816
$wsExhibit.Range("F31").Value = 799
$wsExhibit.Range("F39").Value = 4597

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 13
$wsAll.Range("F10").Value = 7791
$wsAll.Range("F13").Value = 6772
$wsAll.Range("F16").Value = 4808
$wsAll.Range("F17").Value = 5172
$wsAll.Range("F25").Value = 90
$wsAll.Range("F27").Value = 8816
$wsAll.Range("F32").Value = 799
$wsAll.Range("F39").Value = 4597
